# Apply updated crypto market data (price + 1h volume change) scraped on
# Mon Aug 26 09:35:19 UTC 2024. Row 47/48 (Stellar / VeChain) also swap rank order.
#
# Numeric-looking price strings (e.g. "561.37", "1.00") are written with a
# leading apostrophe so Excel keeps them as literal text (quote-prefix), just
# like the source workbook stores them -- otherwise Excel would silently
# coerce them into real numbers and normalise away trailing zeros / thousand dots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "63.558.91"
$ws.Cells.Item(2, 5).Value = "  -0.37%  "

$ws.Cells.Item(3, 4).Value = "2.724.02"
$ws.Cells.Item(3, 5).Value = "  -0.79%  "

$ws.Cells.Item(4, 5).Value = "  +0.09%  "

$ws.Cells.Item(5, 4).Value = "'561.37"
$ws.Cells.Item(5, 5).Value = "  -2.01%  "

$ws.Cells.Item(6, 4).Value = "'158.95"
$ws.Cells.Item(6, 5).Value = "  +1.22%  "

$ws.Cells.Item(7, 5).Value = "  +0.11%  "

$ws.Cells.Item(8, 4).Value = "'0.595"
$ws.Cells.Item(8, 5).Value = "  -0.96%  "

$ws.Cells.Item(9, 5).Value = "  -0.03%  "

$ws.Cells.Item(10, 4).Value = "'0.166"
$ws.Cells.Item(10, 5).Value = "  +4.16%  "

$ws.Cells.Item(11, 4).Value = "'5.61"
$ws.Cells.Item(11, 5).Value = "  +2.19%  "

$ws.Cells.Item(12, 4).Value = "'0.378"
$ws.Cells.Item(12, 5).Value = "  -0.72%  "

$ws.Cells.Item(13, 4).Value = "3.208.11"
$ws.Cells.Item(13, 5).Value = "  -0.72%  "

$ws.Cells.Item(14, 4).Value = "'26.65"
$ws.Cells.Item(14, 5).Value = "  +1.00%  "

$ws.Cells.Item(15, 4).Value = "63.409.72"
$ws.Cells.Item(15, 5).Value = "  -0.13%  "

$ws.Cells.Item(16, 5).Value = "  -0.16%  "

$ws.Cells.Item(17, 4).Value = "2.727.61"
$ws.Cells.Item(17, 5).Value = "  -0.82%  "

$ws.Cells.Item(18, 4).Value = "'12.48"
$ws.Cells.Item(18, 5).Value = "  +2.99%  "

$ws.Cells.Item(19, 4).Value = "'4.71"
$ws.Cells.Item(19, 5).Value = "  -1.64%  "

$ws.Cells.Item(20, 4).Value = "'352.66"
$ws.Cells.Item(20, 5).Value = "  -0.26%  "

$ws.Cells.Item(21, 4).Value = "'6.54"
$ws.Cells.Item(21, 5).Value = "  -2.57%  "

$ws.Cells.Item(22, 5).Value = "  +0.22%  "

$ws.Cells.Item(23, 4).Value = "'0.519"
$ws.Cells.Item(23, 5).Value = "  -2.99%  "

$ws.Cells.Item(24, 4).Value = "'64.22"
$ws.Cells.Item(24, 5).Value = "  -1.27%  "

$ws.Cells.Item(25, 5).Value = "  +0.17%  "

$ws.Cells.Item(26, 4).Value = "'0.999"
$ws.Cells.Item(26, 5).Value = "  -0.05%  "

$ws.Cells.Item(27, 4).Value = "'8.32"
$ws.Cells.Item(27, 5).Value = "  -0.52%  "

$ws.Cells.Item(28, 4).Value = "0.0₃0898"
$ws.Cells.Item(28, 5).Value = "  -0.42%  "

$ws.Cells.Item(29, 5).Value = "  +1.29%  "

$ws.Cells.Item(30, 4).Value = "'7.12"
$ws.Cells.Item(30, 5).Value = "  +2.97%  "

$ws.Cells.Item(32, 4).Value = "'165.95"
$ws.Cells.Item(32, 5).Value = "  -1.89%  "

$ws.Cells.Item(33, 4).Value = "'19.97"
$ws.Cells.Item(33, 5).Value = "  -0.68%  "

$ws.Cells.Item(34, 4).Value = "'4.85"
$ws.Cells.Item(34, 5).Value = "  +0.46%  "

$ws.Cells.Item(35, 4).Value = "'0.999"
$ws.Cells.Item(35, 5).Value = "  +0.01%  "

$ws.Cells.Item(36, 5).Value = "  +2.53%  "

$ws.Cells.Item(37, 5).Value = "  +0.32%  "

$ws.Cells.Item(38, 4).Value = "'0.967"
$ws.Cells.Item(38, 5).Value = "  -0.85%  "

$ws.Cells.Item(39, 4).Value = "'343.11"
$ws.Cells.Item(39, 5).Value = "  +5.71%  "

$ws.Cells.Item(40, 4).Value = "'6.25"
$ws.Cells.Item(40, 5).Value = "  +1.94%  "

$ws.Cells.Item(41, 5).Value = "  -1.34%  "

$ws.Cells.Item(42, 4).Value = "'38.51"
$ws.Cells.Item(42, 5).Value = "  -0.85%  "

$ws.Cells.Item(43, 4).Value = "'21.67"
$ws.Cells.Item(43, 5).Value = "  +2.17%  "

$ws.Cells.Item(44, 4).Value = "'20.96"
$ws.Cells.Item(44, 5).Value = "  -1.16%  "

$ws.Cells.Item(45, 4).Value = "'0.0580"
$ws.Cells.Item(45, 5).Value = "  -0.58%  "

$ws.Cells.Item(46, 4).Value = "'0.624"
$ws.Cells.Item(46, 5).Value = "  +0.41%  "

$ws.Cells.Item(47, 2).Value = "VeChain"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(47, 4).Value = "'0.0249"
$ws.Cells.Item(47, 5).Value = "  -1.45%  "

$ws.Cells.Item(48, 2).Value = "Stellar"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(48, 4).Value = "'0.0995"
$ws.Cells.Item(48, 5).Value = "  -0.79%  "

$ws.Cells.Item(49, 5).Value = "  +0.11%  "

$ws.Cells.Item(51, 4).Value = "'131.06"
$ws.Cells.Item(51, 5).Value = "  -2.61%  "
